$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered to exclude save games), columns B:E and G (sum).
# Column F (Win) is left untouched.

$data = @{
    2 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    3 = @(0.6606524410359556, 10.34677158129881, 3.537761648806719, 10.19245300693656, 24.73763867807805)
    4 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    5 = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 10.19245300693656, 37.52559925287081)
    6 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 10.19245300693656, 13.45301510845117)
    7 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548)
    8 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
